$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148; this shifts existing rows 148:214 down to 149:215
$ws.Rows("148:148").Insert()

# Populate the newly inserted row 148 with the new data values
$ws.Range("A148").Value = 3
$ws.Range("B148").Value = "Femacal de La Calera"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = (Get-Date -Year 2021 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E148").Value = 5
$ws.Range("F148").Value = 100112032
$ws.Range("G148").Value = "Zapallo italiano"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 160
$ws.Range("K148").Value = 10000
$ws.Range("L148").Value = 11000
$ws.Range("M148").Value = 10500
$ws.Range("N148").Value = "`$/caja 70 unidades"
$ws.Range("O148").Value = "Región de Arica y Parinacota"
$ws.Range("P148").Value = 150
$ws.Range("Q148").Value = 70
$ws.Range("R148").Value = "Hortaliza"

# Copy the style (number format) of the date cell from the row below, to match other date cells
$ws.Range("D149").Copy()
$ws.Range("D148").PasteSpecial(-4122)  # xlPasteFormats
